$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value  = "0 - 4"
$ws.Range("A3").Value  = "5 - 9"
$ws.Range("A4").Value  = "10 - 14"
$ws.Range("A5").Value  = "15 - 19"
$ws.Range("A6").Value  = "20 - 24"
$ws.Range("A7").Value  = "25 - 29"
$ws.Range("A8").Value  = "30 - 34"
$ws.Range("A9").Value  = "35 - 39"
$ws.Range("A10").Value = "40 - 44"
$ws.Range("A11").Value = "45 - 49"
$ws.Range("A12").Value = "50 - 54"
$ws.Range("A13").Value = "55 - 59"
$ws.Range("A14").Value = "60 - 64"
$ws.Range("A15").Value = "65 - 69"
$ws.Range("A16").Value = "70 - 74"
$ws.Range("A17").Value = "75 - 79"
$ws.Range("A18").Value = "80 - 84"
$ws.Range("A19").Value = "85 - 89"
$ws.Range("A20").Value = "90 - 94"
$ws.Range("A21").Value = "95 - 99"
$ws.Range("A22").Value = "100+"
